# ADD results from server
# Updates the header row (row 1) and the data row (row 2) on each of the
# six yearly worksheets (2025, 2030, 2035, 2040, 2045, 2050) with a fresh
# set of results: two new technology columns ("gb" and "btes") are
# introduced while two old ones ("gt" and "dgt") are dropped, and all the
# investment-cost figures in row 2 are replaced with newly computed values.

$wb = $excel.ActiveWorkbook

# New column headers (A1:O1) shared by every sheet.
$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# New row-2 values (A2:O2) per sheet, in sheet order (2025, 2030, 2035, 2040, 2045, 2050).
$sheetValues = @(
    @(39063.99109145206, 0, 483537.6274462014, 0, 2897240.114301849, 94331.34471502228, 0, 25342.77928792104, 0, 0, 0, 0, 0, 23638.06126801545, 19940.13531829346),
    @(30846.52922536713, 0, 1495599.874611417, 0, 0, 70193.79982138964, 0, 56602.42752520426, 0, 0, 0, 0, 0, 51649.16401227913, 42574.77934331147),
    @(242452.4252219552, 0, 943335.270081223, 0, 0, 1425.925979620855, 0, 39373.98526588717, 0, 0, 0, 0, 0, 53308.16490721726, 30023.09380555204),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11578.49752443177, 0),
    @(76705.58894163162, 1930.947398408091, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 28147.3462746636, 8312.661449003012),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $values = $sheetValues[$s - 1]

    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
        $ws.Cells.Item(2, $i + 1).Value = $values[$i]
    }
}

Write-Output "Updated header row and data row on $($wb.Worksheets.Count) sheets."
